# optimizing and new programs 10feb
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Done [yes or no]" status column (C) for rows 21-41 with new values,
# in top-to-bottom order so new shared-string entries are created in the
# same order they appear in the target workbook.
$ws.Range("C21").Value = "not done"
$ws.Range("C22").Value = " done"

$ws.Range("C27").Value = "not done"
$ws.Range("C28").Value = "done"
$ws.Range("C29").Value = "done"
$ws.Range("C30").Value = "done"
$ws.Range("C31").Value = "done"
$ws.Range("C32").Value = "done"
$ws.Range("C33").Value = "done"
$ws.Range("C34").Value = "Done"
$ws.Range("C35").Value = "DONE"
$ws.Range("C36").Value = "Done"
$ws.Range("C37").Value = "done"
$ws.Range("C38").Value = "Done"
$ws.Range("C39").Value = "Done"
$ws.Range("C40").Value = "Done"
$ws.Range("C41").Value = "need to be optimized"

# Rows 42-43 no longer carry a status value.
$ws.Range("C42").ClearContents()
$ws.Range("C43").ClearContents()

# Update the view: scroll position / selected cell.
$ws.Range("B35").Select()
